# This script re-orders the data rows (A2:M24) of the active worksheet.
# The original export grouped the records by "Repository" (K column) and
# then by "Version"/timestamp (L/M columns) within each group. The updated
# analytics data (regenerated while adding the graph labels / total graph)
# instead lists every record sorted strictly by the "Version" timestamp,
# across all repositories - exactly matching the order produced by
# pandas' `DataFrame.sort_values("Version")` (default unstable quicksort).
#
# Rather than trying to reproduce that exact unstable-sort tie-break with
# Excel's own Range.Sort (which uses a different, stable algorithm and
# would not reproduce the same tie-break ordering), we simply write out
# the final, already-sorted values for every row/column directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.625,0,1,1,1,1,0.268125,0.5,0.7681249999999999,70,"User","2021-09-14-03:28","2021-09-14-03:28"),
    @(0.7142857142857143,0,0.8928571428571429,0.8928571428571429,0.8928571428571429,0.8928571428571429,0.2651785714285714,0.4464285714285714,0.7116071428571429,629,"Content","2021-09-14-03:28","2021-09-14-03:28"),
    @(0.5555555555555556,0,1,1,1,1,0.2566666666666667,0.5,0.7566666666666667,50,"Mobile_App","2021-09-14-03:28","2021-09-14-03:28"),
    @(0.5714285714285714,0,1,1,1,1,0.2592857142857143,0.5,0.7592857142857143,54,"Files","2021-09-14-03:28","2021-09-14-03:28"),
    @(0.625,0,1,1,1,1,0.268125,0.5,0.7681249999999999,70,"User","2021-09-27-02.31","2021-09-27-02.31"),
    @(0.5555555555555556,0,1,1,1,1,0.2566666666666667,0.5,0.7566666666666667,50,"Mobile_App","2021-09-27-03.13","2021-09-27-03.13"),
    @(0.7115384615384616,0,0.9423076923076923,0.9423076923076923,0.9423076923076923,0.9423076923076923,0.2728846153846154,0.4711538461538461,0.7440384615384615,1323,"Content","2021-09-27-03.22","2021-09-27-03.22"),
    @(0.7230769230769231,0,0.9692307692307692,0.9692307692307692,0.9692307692307692,0.9692307692307692,0.2792307692307692,0.4846153846153846,0.7638461538461538,1879,"Content","2021-10-06-19.24","2021-10-06-19.24"),
    @(0.5,0,1,1,1,1,0.2475,0.5,0.7475000000000001,125,"Files","2021-10-15-19.42","2021-10-15-19.42"),
    @(0.5,0,1,1,1,1,0.2475,0.5,0.7475000000000001,125,"Files","2021-10-15-20.05","2021-10-15-20.05"),
    @(0.5185185185185185,0,1,1,1,1,0.2505555555555555,0.5,0.7505555555555555,322,"Admin","2021-10-18-01.33","2021-10-18-01.33"),
    @(0.5411764705882353,0,1,1,1,1,0.2542941176470588,0.5,0.7542941176470588,2342,"Mobile_App","2021-10-19-18.47","2021-10-19-18.47"),
    @(0.6979166666666666,0,0.9791666666666666,0.9791666666666666,0.9791666666666666,0.9791666666666666,0.27671875,0.4895833333333333,0.7663020833333333,2477,"Content","2021-10-19-18.48","2021-10-19-18.48"),
    @(0.5,0,0.9285714285714286,0.9285714285714286,0.9285714285714286,0.9285714285714286,0.2357142857142857,0.4642857142857143,0.7,685,"User","2021-10-24-19.18","2021-10-24-19.18"),
    @(0.5411764705882353,0,1,1,1,1,0.2542941176470588,0.5,0.7542941176470588,2222,"Mobile_App","2021-11-03-22:59","2021-11-03-22:59"),
    @(0.5,0,1,1,1,1,0.2475,0.5,0.7475000000000001,132,"Files","2021-11-03-22:59","2021-11-03-22:59"),
    @(0.71,0,0.98,0.98,0.98,0.98,0.27885,0.49,0.76885,2647,"Content","2021-11-03-22:59","2021-11-03-22:59"),
    @(0.5,0,0.9285714285714286,0.9285714285714286,0.9285714285714286,0.9285714285714286,0.2357142857142857,0.4642857142857143,0.7,688,"User","2021-11-03-22:59","2021-11-03-22:59"),
    @(0.6744186046511628,0,1,1,1,1,0.2762790697674419,0.5,0.7762790697674419,708,"Admin","2021-11-04-00.35","2021-11-04-00.35"),
    @(0.7037037037037037,0,1,1,1,1,0.2811111111111111,0.5,0.7811111111111111,1051,"Admin","2021-11-08-03.34","2021-11-08-03.34"),
    @(0.5402298850574713,0,1,1,1,1,0.2541379310344828,0.5,0.7541379310344828,2340,"Mobile_App","2021-11-08-20.15","2021-11-08-20.15"),
    @(0.7032967032967034,0.02197802197802198,0.978021978021978,0.978021978021978,0.978021978021978,0.978021978021978,0.281043956043956,0.4890109890109889,0.770054945054945,2038,"Content","2021-11-10-17.00","2021-11-10-17.00"),
    @(0.5444444444444444,0,1,1,1,1,0.2548333333333334,0.5,0.7548333333333334,2439,"Mobile_App","2021-11-10-17.28","2021-11-10-17.28")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowValues[$j]
    }
}
